$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.977.94'
$ws.Range('E2').Value = '  +2.66%  '
$ws.Range('D3').Value = '2.249.66'
$ws.Range('E3').Value = '  +1.85%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'269.80"
$ws.Range('E5').Value = '  +5.37%  '
$ws.Range('D6').Value = '87.85'
$ws.Range('E6').Value = '  +13.83%  '
$ws.Range('D7').Value = '0.618'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +4.11%  '
$ws.Range('D10').Value = '45.97'
$ws.Range('E10').Value = '  +7.76%  '
$ws.Range('E11').Value = '  +2.32%  '
$ws.Range('D12').Value = '7.61'
$ws.Range('E12').Value = '  +9.22%  '
$ws.Range('E13').Value = '  +2.02%  '
$ws.Range('D14').Value = '2.581.34'
$ws.Range('E14').Value = '  +1.48%  '
$ws.Range('D15').Value = '15.06'
$ws.Range('E15').Value = '  +4.70%  '
$ws.Range('D16').Value = '2.236.95'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('D17').Value = '0.801'
$ws.Range('E17').Value = '  +2.42%  '
$ws.Range('D18').Value = '43.970.52'
$ws.Range('E18').Value = '  +2.86%  '
$ws.Range('E19').Value = '  +0.76%  '
$ws.Range('E20').Value = '  +1.72%  '
$ws.Range('D21').Value = '70.25'
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('E22').Value = '  +4.18%  '
$ws.Range('D23').Value = "'234.00"
$ws.Range('E23').Value = '  +1.92%  '
$ws.Range('D24').Value = '8.82'
$ws.Range('E24').Value = '  -4.04%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').Value = '2.56'
$ws.Range('E26').Value = '  +16.59%  '
$ws.Range('E27').Value = '  +2.85%  '
$ws.Range('E28').Value = '  +6.30%  '
$ws.Range('D29').Value = '40.71'
$ws.Range('E29').Value = '  -4.18%  '
$ws.Range('D30').Value = '2.32'
$ws.Range('E30').Value = '  +5.36%  '
$ws.Range('D31').Value = '175.19'
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('D32').Value = '0.0912'
$ws.Range('E32').Value = '  +4.19%  '
$ws.Range('E33').Value = '  +2.40%  '
$ws.Range('D34').Value = '5.41'
$ws.Range('E34').Value = '  +4.29%  '
$ws.Range('D35').Value = '0.124'
$ws.Range('E35').Value = '  +2.34%  '
$ws.Range('E36').Value = '  +5.33%  '
$ws.Range('D37').Value = '0.0356'
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('D38').Value = '4.43'
$ws.Range('E38').Value = '  +1.76%  '
$ws.Range('E39').Value = '  +18.21%  '
$ws.Range('D40').Value = '12.74'
$ws.Range('E40').Value = '  -2.44%  '
$ws.Range('E41').Value = '  +2.51%  '
$ws.Range('D42').Value = '65.35'
$ws.Range('E42').Value = '  +7.51%  '
$ws.Range('E43').Value = '  +2.53%  '
$ws.Range('E44').Value = '  +2.28%  '
$ws.Range('D45').Value = "'0.0990"
$ws.Range('E45').Value = '  +2.12%  '
$ws.Range('D46').Value = "'8.40"
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('D47').Value = '100.51'
$ws.Range('E47').Value = '  -2.19%  '
$ws.Range('D48').Value = '1.22'
$ws.Range('E48').Value = '  +7.50%  '
$ws.Range('E49').Value = '  +2.57%  '
$ws.Range('D50').Value = '1.53'
$ws.Range('E50').Value = '  +2.87%  '
$ws.Range('E51').Value = '  -5.51%  '
